$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.869.20"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "1.638.63"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("E4").Value = "  +0.43%  "
$ws.Range("D5").Value = "215.44"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.59%  "
$ws.Range("D8").Value = "28.76"
$ws.Range("E8").Value = "  -1.99%  "
$ws.Range("E9").Value = "  +0.69%  "
$ws.Range("D10").Value = "0.0608"
$ws.Range("E10").Value = "  +0.33%  "
$ws.Range("E11").Value = "  -1.09%  "
$ws.Range("D12").Value = "1.874.43"
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("D13").Value = "1.631.91"
$ws.Range("E13").Value = "  +0.29%  "
$ws.Range("E14").Value = "  +4.53%  "
$ws.Range("D15").Value = "9.46"
$ws.Range("E15").Value = "  +7.66%  "
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("D17").Value = "29.877.97"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "64.60"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").Value = "240.36"
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("D22").Value = "9.89"
$ws.Range("E22").Value = "  +3.37%  "
$ws.Range("D23").Value = "4.15"
$ws.Range("E23").Value = "  +1.11%  "
$ws.Range("D24").Value = "2.20"
$ws.Range("E24").Value = "  +2.94%  "
$ws.Range("D25").Value = "157.84"
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("D26").Value = "15.52"
$ws.Range("E26").Value = "  -0.43%  "
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("E28").Value = "  +0.68%  "
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("E32").Value = "  +1.93%  "
$ws.Range("D33").Value = "3.19"
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("D34").Value = "1.425.10"
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("E35").Value = "  +3.69%  "
$ws.Range("E36").Value = "  -1.02%  "
$ws.Range("E37").Value = "  -4.43%  "
$ws.Range("E38").Value = "  +2.70%  "
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("D40").Value = "76.61"
$ws.Range("E40").Value = "  +10.97%  "
$ws.Range("D41").Value = "0.560"
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("D42").Value = "0.834"
$ws.Range("E42").Value = "  +1.24%  "
$ws.Range("E43").Value = "  -1.35%  "
$ws.Range("E44").Value = "  -0.28%  "
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("E46").Value = "  -1.59%  "
$ws.Range("E47").Value = "  -0.61%  "
$ws.Range("D48").Value = "1.781.65"
$ws.Range("E48").Value = "  +0.96%  "
$ws.Range("D49").Value = "48.91"
$ws.Range("E49").Value = "  -9.19%  "
$ws.Range("D50").Value = "93.53"
$ws.Range("E50").Value = "  +6.25%  "
$ws.Range("E51").Value = "  -1.37%  "
